$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "Table 1" sheet: fix Yes -> no for I193:I197, then append 60 new
#    contaminant rows (360-419) mirroring the pattern of the existing
#    rows immediately above them.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1")

$ws1.Range("I193:I197").Value = "no"

$newNames = @(
    '1,2,3-Trichloropropane',
    'Chromium (hexavalent)',
    'Perfluorinated chemicals',
    'Total trihalomethanes (TTHMs)',
    '1,4-Dioxane',
    'Bromodichloromethane',
    'Chloroform',
    'Dichloroacetic acid',
    'Trichloroacetic acid',
    'Radiological contaminants',
    'Manganese',
    'Dibromochloromethane',
    'Tetrachloroethylene (perchloroethylene)',
    'Haloacetic acids (HAA5)',
    'Hormones',
    'Arsenic',
    'Trichloroethylene',
    'Atrazine',
    'Chlorate',
    'Thallium',
    'Chlorite',
    'Bromate',
    'Bromoform',
    'Nitrate',
    'Benzo[a]pyrene',
    'Nitrate and nitrite',
    'Cadmium',
    'Beryllium',
    'Antimony',
    'Strontium',
    'Benzene',
    'Mercury (inorganic)',
    'Perchlorate',
    'Carbon tetrachloride',
    'Vinyl chloride',
    '1,2-Dichloroethane',
    'Barium',
    'Selenium',
    'Chloromethane',
    'Heptachlor epoxide',
    'Fluoride',
    'Aluminum',
    '1,1,2-Trichloroethane',
    'Chlordane',
    'Lindane',
    'Heptachlor',
    '1,2-Dichloropropane',
    'Styrene',
    '1,2-Dibromo-3-chloropropane (DBCP)',
    'Ethylene dibromide',
    'Molybdenum',
    'Polychlorinated biphenyls (PCBs)',
    '1,3-Butadiene',
    'Di(2-ethylhexyl) phthalate',
    'Carbofuran',
    'p-Dichlorobenzene',
    'Glyphosate',
    'Vanadium',
    'MTBE',
    'N-Nitrosodimethylamine (NDMA)'
)

$startRow = 360
for ($i = 0; $i -lt $newNames.Count; $i++) {
    $r = $startRow + $i
    $ws1.Cells.Item($r, 2).Value = $newNames[$i]
    $ws1.Cells.Item($r, 5).Value = "water pollution"
    $ws1.Cells.Item($r, 6).Value = "physical environment"
    $ws1.Cells.Item($r, 9).Value = "Yes"
}

$endRow = $startRow + $newNames.Count - 1

# Match the styling used by the existing table rows:
#  - column B uses the plain "Normal" style
#  - column F uses the same style as the neighbouring "physical
#    environment" cells in column E/F higher up the sheet (style index 3)
$bRange = $ws1.Range($ws1.Cells.Item($startRow, 2), $ws1.Cells.Item($endRow, 2))
$bRange.Style = "Normal"

$fmtSrc = $ws1.Cells.Item($startRow - 1, 6)
$fRange = $ws1.Range($ws1.Cells.Item($startRow, 6), $ws1.Cells.Item($endRow, 6))
$fmtSrc.Copy()
$fRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A193").Select()

# ------------------------------------------------------------------
# 2) "Sheet2" sheet: replace the old H13:H39 lookup list with the new
#    F5:BM66 layout - a 60-column header row of the same 60 new names,
#    followed by the same 60 names listed down column F.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Cells.Clear()

# Drop the bestFit/customWidth override that used to live on column I
# (index 9) so only column H (index 8) keeps its custom width, matching
# the target layout.
$ws2.Columns.Item(9).Delete()
$ws2.Columns.Item(9).Insert()

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $col = 6 + $i
    $ws2.Cells.Item(5, $col).Value = $newNames[$i]
}

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $row = 7 + $i
    $ws2.Cells.Item($row, 6).Value = $newNames[$i]
}

$ws2.Range("F5:BM66").Style = "Normal"

$ws2.Range("F7").Select()
